$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix two typo'd tag values in column F (",(E-7H)" -> "RT,(E-7H)") ---
$ws.Range("F111").Value = "RT,(E-7H)"
$ws.Range("F117").Value = "RT,(E-7H)"

# --- Append a new product row (row 124) for group 21, matching the
#     layout/formatting of the row above it (row 123) ---
$ws.Range("A123:F123").Copy() | Out-Null
$ws.Range("A124:F124").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Columns A, D and E hold numeric-looking codes that must stay text
# (same as every other cell in this table), so enter them with a
# leading apostrophe to force text entry instead of letting Excel
# auto-convert them to numbers.
$ws.Range("A124").Formula = "'20141057"
$ws.Range("B124").Value = "KZLER SS TOM YUM 60G"
$ws.Range("C124").Value = "CLP03N"
$ws.Range("D124").Formula = "'21"
$ws.Range("E124").Formula = "'3"
$ws.Range("F124").Value = "RT,(E-7H)"
